# Sweden_M2.xlsx update:
#  - Row 307 (2023-06-01): revise O/H/L/C values from 4811874000000 to 4812070000000
#  - Append three new monthly rows (308-310) for 2023-07-01, 2023-08-01, 2023-09-01
#    continuing the "ECONOMICS:SEM2" M2 series, with volume 0.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 307 (open/high/low/close columns) ---
$ws.Range("C307").Value = 4812070000000
$ws.Range("D307").Value = 4812070000000
$ws.Range("E307").Value = 4812070000000
$ws.Range("F307").Value = 4812070000000

# --- Carry the row 307 number formatting / styling down onto the new rows ---
$ws.Range("A307:G307").Copy()
$ws.Range("A308:G310").PasteSpecial(-4122)

# --- Row 308: 2023-07-01 ---
$ws.Range("A308").Value = 45108.41666666666
$ws.Range("B308").Value = "ECONOMICS:SEM2"
$ws.Range("C308").Value = 4784709000000
$ws.Range("D308").Value = 4784709000000
$ws.Range("E308").Value = 4784709000000
$ws.Range("F308").Value = 4784709000000
$ws.Range("G308").Value = 0

# --- Row 309: 2023-08-01 ---
$ws.Range("A309").Value = 45139.41666666666
$ws.Range("B309").Value = "ECONOMICS:SEM2"
$ws.Range("C309").Value = 4765034000000
$ws.Range("D309").Value = 4765034000000
$ws.Range("E309").Value = 4765034000000
$ws.Range("F309").Value = 4765034000000
$ws.Range("G309").Value = 0

# --- Row 310: 2023-09-01 ---
$ws.Range("A310").Value = 45170.41666666666
$ws.Range("B310").Value = "ECONOMICS:SEM2"
$ws.Range("C310").Value = 4744031000000
$ws.Range("D310").Value = 4744031000000
$ws.Range("E310").Value = 4744031000000
$ws.Range("F310").Value = 4744031000000
$ws.Range("G310").Value = 0
